# Scheduled-runner style data refresh for the Gungnir_Profits workbook.
# Updates cached market/profit figures (columns H:N) on several leve rows
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 607.5
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 619.44446
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 619.44446
$ws.Range("M19").Value = -325
$ws.Range("N19").Value = -969.44446

$ws.Range("H28").Value = 300.11765
$ws.Range("I28").Value = 260.13333
$ws.Range("J28").Value = 600
$ws.Range("K28").Value = 260.13333
$ws.Range("L28").Value = 600
$ws.Range("M28").Value = 224.86667
$ws.Range("N28").Value = -1570

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H80").Value = 319.16666
$ws.Range("I80").Value = 315.77777
$ws.Range("J80").Value = 322.55554
$ws.Range("K80").Value = 947.33331
$ws.Range("L80").Value = 967.66662
$ws.Range("M80").Value = 50.66669000000002
$ws.Range("N80").Value = -2963.66662

$ws.Range("H83").Value = 319.16666
$ws.Range("I83").Value = 315.77777
$ws.Range("J83").Value = 322.55554
$ws.Range("K83").Value = 2841.99993
$ws.Range("L83").Value = 2902.99986
$ws.Range("M83").Value = 2150.00007
$ws.Range("N83").Value = -12886.99986

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31256170
$ws.Range("I32").Value = 17246896
$ws.Range("J32").Value = 166679170
$ws.Range("K32").Value = 17246896
$ws.Range("L32").Value = 166679170
$ws.Range("M32").Value = -17246609
$ws.Range("N32").Value = -166679744

$ws.Range("H45").Value = 43032.168
$ws.Range("I45").Value = 101207.4
$ws.Range("J45").Value = 1478.4286
$ws.Range("K45").Value = 101207.4
$ws.Range("L45").Value = 1478.4286
$ws.Range("M45").Value = -100830.4
$ws.Range("N45").Value = -2232.4286

$ws.Range("H102").Value = 1475
$ws.Range("I102").Value = 1475
$ws.Range("K102").Value = 1475
$ws.Range("M102").Value = 147

$ws.Range("H122").Value = 721.0476
$ws.Range("I122").Value = 767.4286
$ws.Range("J122").Value = 628.2857
$ws.Range("K122").Value = 2302.2858
$ws.Range("L122").Value = 1884.8571
$ws.Range("M122").Value = 147.7142000000003
$ws.Range("N122").Value = -6784.8571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1036
$ws.Range("I94").Value = 1039.6
$ws.Range("K94").Value = 1039.6
$ws.Range("M94").Value = -588.5999999999999

$ws.Range("H99").Value = 1196.6666
$ws.Range("I99").Value = 1055
$ws.Range("J99").Value = 1237.1428
$ws.Range("K99").Value = 1055
$ws.Range("L99").Value = 1237.1428
$ws.Range("M99").Value = 443
$ws.Range("N99").Value = -4233.1428

$ws.Range("H105").Value = 111113260
$ws.Range("I105").Value = 2100
$ws.Range("J105").Value = 250002220
$ws.Range("K105").Value = 2100
$ws.Range("L105").Value = 250002220
$ws.Range("M105").Value = -353
$ws.Range("N105").Value = -250005714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3741.3572
$ws.Range("I62").Value = 2816.5
$ws.Range("J62").Value = 4435
$ws.Range("K62").Value = 2816.5
$ws.Range("L62").Value = 4435
$ws.Range("M62").Value = -2192.5
$ws.Range("N62").Value = -5683

$ws.Range("H65").Value = 3741.3572
$ws.Range("I65").Value = 2816.5
$ws.Range("J65").Value = 4435
$ws.Range("K65").Value = 14082.5
$ws.Range("L65").Value = 22175
$ws.Range("M65").Value = -10962.5
$ws.Range("N65").Value = -28415

$ws.Range("H105").Value = 54990
$ws.Range("I105").Value = 100000
$ws.Range("J105").Value = 9980
$ws.Range("K105").Value = 100000
$ws.Range("L105").Value = 9980
$ws.Range("M105").Value = -98253
$ws.Range("N105").Value = -13474

$ws.Range("H107").Value = 891.2
$ws.Range("I107").Value = 752
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 752
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 1168
$ws.Range("N107").Value = -4940

$ws.Range("H122").Value = 31250902
$ws.Range("I122").Value = 31250902
$ws.Range("K122").Value = 93752706
$ws.Range("M122").Value = -93750256

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 41029804
$ws.Range("I107").Value = 160.5
$ws.Range("J107").Value = 66673332
$ws.Range("K107").Value = 481.5
$ws.Range("L107").Value = 200019996
$ws.Range("M107").Value = 1438.5
$ws.Range("N107").Value = -200023836

$ws.Range("H132").Value = 17246182
$ws.Range("I132").Value = 962.2143
$ws.Range("J132").Value = 33341718
$ws.Range("K132").Value = 8659.9287
$ws.Range("L132").Value = 300075462
$ws.Range("M132").Value = -6129.9287
$ws.Range("N132").Value = -300080522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8622116
$ws.Range("I40").Value = 1502.3
$ws.Range("J40").Value = 27779034
$ws.Range("K40").Value = 1502.3
$ws.Range("L40").Value = 27779034
$ws.Range("M40").Value = -1366.3
$ws.Range("N40").Value = -27779306

$ws.Range("H68").Value = 1275.2693
$ws.Range("I68").Value = 1087.2106
$ws.Range("K68").Value = 1087.2106
$ws.Range("M68").Value = -338.2106000000001

$ws.Range("H71").Value = 1275.2693
$ws.Range("I71").Value = 1087.2106
$ws.Range("K71").Value = 5436.053000000001
$ws.Range("M71").Value = -1692.053000000001

$ws.Range("H136").Value = 40819576
$ws.Range("I136").Value = 5497929.5
$ws.Range("J136").Value = 500001000
$ws.Range("K136").Value = 16493788.5
$ws.Range("L136").Value = 1500003000
$ws.Range("M136").Value = -16491238.5
$ws.Range("N136").Value = -1500008100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1280.3
$ws.Range("I81").Value = 1250.5
$ws.Range("J81").Value = 1325
$ws.Range("K81").Value = 2501
$ws.Range("L81").Value = 2650
$ws.Range("M81").Value = -1440
$ws.Range("N81").Value = -4772

$ws.Range("H84").Value = 1280.3
$ws.Range("I84").Value = 1250.5
$ws.Range("J84").Value = 1325
$ws.Range("K84").Value = 12505
$ws.Range("L84").Value = 13250
$ws.Range("M84").Value = -7201
$ws.Range("N84").Value = -23858

$ws.Range("H122").Value = 23289.2
$ws.Range("I122").Value = 56510.89
$ws.Range("K122").Value = 169532.67
$ws.Range("M122").Value = -167082.67

$ws.Range("H140").Value = 42841.3
$ws.Range("J140").Value = 42841.3
$ws.Range("L140").Value = 42841.3
$ws.Range("N140").Value = -53201.3
